$d = $word.ActiveDocument

function Repeat-Space($n) {
    $s = ""
    for ($i = 0; $i -lt $n; $i++) {
        $s = $s + " "
    }
    return $s
}

function Set-ParaXml($paraIndex, $innerXml) {
    $para = $d.Paragraphs.Item($paraIndex)
    $r = $para.Range
    $pkg = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $innerXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($pkg)
}

# --- 1. PHARMA SALES DASHBOARD heading: merge the three runs into one ---
$pharmaText = (Repeat-Space 58) + "PHARMA SALES DASHBOARD"
$pharmaXml = '<w:body><w:p><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">' + $pharmaText + '</w:t></w:r></w:p></w:body>'
Set-ParaXml 14 $pharmaXml

Write-Output ("Para14: [" + $d.Paragraphs.Item(14).Range.Text + "]")
